$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the 3 extra rows (8-10) that no longer exist in the updated data
$ws.Range("A8:T10").Delete()

# Update rows 2-7 with the refreshed TPM-derived values
$ws.Cells.Item(2,1).Value = "FAPs"
$ws.Cells.Item(2,2).Value = "Angpt4"
$ws.Cells.Item(2,3).Value = "Tek"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 2.365546333333333
$ws.Cells.Item(2,8).Value = 7.096639
$ws.Cells.Item(2,9).Value = 0.6866153222176684
$ws.Cells.Item(2,10).Value = 0.6866153222176683
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 53.457377
$ws.Cells.Item(2,14).Value = 160.372131
$ws.Cells.Item(2,15).Value = 0.6217639481372091
$ws.Cells.Item(2,16).Value = 0.6217639481372091
$ws.Cells.Item(2,17).Value = 126.4559021519676
$ws.Cells.Item(2,18).Value = 1138.103119367709
$ws.Cells.Item(2,19).Value = 0.4269126535935595
$ws.Cells.Item(2,20).Value = 0.4269126535935595
$ws.Cells.Item(3,1).Value = "FAPs"
$ws.Cells.Item(3,2).Value = "Angpt4"
$ws.Cells.Item(3,3).Value = "Tek"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 2.365546333333333
$ws.Cells.Item(3,8).Value = 7.096639
$ws.Cells.Item(3,9).Value = 0.6866153222176684
$ws.Cells.Item(3,10).Value = 0.6866153222176683
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 31.083557
$ws.Cells.Item(3,14).Value = 93.25067100000001
$ws.Cells.Item(3,15).Value = 0.3615335470438062
$ws.Cells.Item(3,16).Value = 0.3615335470438062
$ws.Cells.Item(3,17).Value = 73.52959428830766
$ws.Cells.Item(3,18).Value = 661.7663485947691
$ws.Cells.Item(3,19).Value = 0.2482344728959796
$ws.Cells.Item(3,20).Value = 0.2482344728959796
$ws.Cells.Item(4,1).Value = "FAPs"
$ws.Cells.Item(4,2).Value = "Angpt4"
$ws.Cells.Item(4,3).Value = "Tek"
$ws.Cells.Item(4,4).Value = "MuSCs"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 2.365546333333333
$ws.Cells.Item(4,8).Value = 7.096639
$ws.Cells.Item(4,9).Value = 0.6866153222176684
$ws.Cells.Item(4,10).Value = 0.6866153222176683
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 1.436030666666667
$ws.Cells.Item(4,14).Value = 4.308092
$ws.Cells.Item(4,15).Value = 0.01670250481898457
$ws.Cells.Item(4,16).Value = 0.01670250481898457
$ws.Cells.Item(4,17).Value = 3.396997078087555
$ws.Cells.Item(4,18).Value = 30.572973702788
$ws.Cells.Item(4,19).Value = 0.01146819572812925
$ws.Cells.Item(4,20).Value = 0.01146819572812925
$ws.Cells.Item(5,1).Value = "MuSCs"
$ws.Cells.Item(5,2).Value = "Angpt4"
$ws.Cells.Item(5,3).Value = "Tek"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 1.079681666666667
$ws.Cells.Item(5,8).Value = 3.239045
$ws.Cells.Item(5,9).Value = 0.3133846777823316
$ws.Cells.Item(5,10).Value = 0.3133846777823315
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 53.457377
$ws.Cells.Item(5,14).Value = 160.372131
$ws.Cells.Item(5,15).Value = 0.6217639481372091
$ws.Cells.Item(5,16).Value = 0.6217639481372091
$ws.Cells.Item(5,17).Value = 57.71694989498833
$ws.Cells.Item(5,18).Value = 519.452549054895
$ws.Cells.Item(5,19).Value = 0.1948512945436496
$ws.Cells.Item(5,20).Value = 0.1948512945436496
$ws.Cells.Item(6,1).Value = "MuSCs"
$ws.Cells.Item(6,2).Value = "Angpt4"
$ws.Cells.Item(6,3).Value = "Tek"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 1.079681666666667
$ws.Cells.Item(6,8).Value = 3.239045
$ws.Cells.Item(6,9).Value = 0.3133846777823316
$ws.Cells.Item(6,10).Value = 0.3133846777823315
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 31.083557
$ws.Cells.Item(6,14).Value = 93.25067100000001
$ws.Cells.Item(6,15).Value = 0.3615335470438062
$ws.Cells.Item(6,16).Value = 0.3615335470438062
$ws.Cells.Item(6,17).Value = 33.56034662768833
$ws.Cells.Item(6,18).Value = 302.043119649195
$ws.Cells.Item(6,19).Value = 0.1132990741478266
$ws.Cells.Item(6,20).Value = 0.1132990741478266
$ws.Cells.Item(7,1).Value = "MuSCs"
$ws.Cells.Item(7,2).Value = "Angpt4"
$ws.Cells.Item(7,3).Value = "Tek"
$ws.Cells.Item(7,4).Value = "MuSCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 1.079681666666667
$ws.Cells.Item(7,8).Value = 3.239045
$ws.Cells.Item(7,9).Value = 0.3133846777823316
$ws.Cells.Item(7,10).Value = 0.3133846777823315
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 1.436030666666667
$ws.Cells.Item(7,14).Value = 4.308092
$ws.Cells.Item(7,15).Value = 0.01670250481898457
$ws.Cells.Item(7,16).Value = 0.01670250481898457
$ws.Cells.Item(7,17).Value = 1.550455983571111
$ws.Cells.Item(7,18).Value = 13.95410385214
$ws.Cells.Item(7,19).Value = 0.005234309090855322
$ws.Cells.Item(7,20).Value = 0.00523430909085532
